# "separate dept from affiliations"
#
# PI hours:
#   - add a new "app" column (F) holding the original (list-of-depts) string
#     that used to live in the "dept" column (E)
#   - "dept" column (E) now holds just the single primary department
# dept hours sheet:
#   - renamed to "department hours"
#   - re-computed: hours/percentage grouped by the single primary dept
#     (so the CSL row disappears and AE's hours change)
# a brand-new sheet "unit(accumulative) hours" is appended at the end,
#   holding exactly what the old "dept hours" sheet used to contain
#   (cumulative hours per department, counting a PI once per every
#   department/affiliation they are listed under)

$wb = $excel.ActiveWorkbook

$wsPI   = $wb.Worksheets.Item(1)
$wsDept = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1) Create the new "unit(accumulative) hours" sheet at the very end,
#    seeded with a copy of the current "dept hours" data (this is the
#    old cumulative-per-department computation, unchanged).
# ------------------------------------------------------------------
$wsAccum = $wb.Worksheets.Add($null, $wsDept)
$wsAccum.Name = "unit(accumulative) hours"

$wsDept.Range("B1:D1").Copy($wsAccum.Range("B1"))
$wsDept.Range("A2:D6").Copy($wsAccum.Range("A2"))
$wsAccum.Range("B1").Value = "unit(accumulative)"

# ------------------------------------------------------------------
# 2) Rewrite "dept hours" -> "department hours" with the new
#    single-primary-dept grouping (ME/AE/ECE/ABE, no CSL row).
# ------------------------------------------------------------------
$wsDept.Name = "department hours"

$wsDept.Cells.Item(2,2).Value = "ME"
$wsDept.Cells.Item(2,3).Value = 62.5
$wsDept.Cells.Item(2,4).Value = 86.20689655172414

$wsDept.Cells.Item(3,2).Value = "AE"
$wsDept.Cells.Item(3,3).Value = 4.5
$wsDept.Cells.Item(3,4).Value = 6.206896551724138

$wsDept.Cells.Item(4,2).Value = "ECE"
$wsDept.Cells.Item(4,3).Value = 3
$wsDept.Cells.Item(4,4).Value = 4.137931034482759

$wsDept.Cells.Item(5,2).Value = "ABE"
$wsDept.Cells.Item(5,3).Value = 2.5
$wsDept.Cells.Item(5,4).Value = 3.448275862068965

# old sheet had 6 rows (incl. CSL), new one only has 5 -> drop the
# trailing row entirely so the sheet dimension shrinks back down
$wsDept.Rows("6:6").Delete()

# ------------------------------------------------------------------
# 3) "PI hours": introduce the "app" column (F) holding the original
#    full affiliation-list strings, and shrink "dept" (E) down to the
#    single primary department for each row.
# ------------------------------------------------------------------
$wsPI.Range("E1").Copy($wsPI.Range("F1"))
$wsPI.Range("F1").Value = "app"

$wsPI.Range("E2:E7").Copy($wsPI.Range("F2:F7"))

$wsPI.Cells.Item(2,5).Value = "ME"
$wsPI.Cells.Item(3,5).Value = "ME"
$wsPI.Cells.Item(4,5).Value = "ME"
$wsPI.Cells.Item(5,5).Value = "AE"
$wsPI.Cells.Item(6,5).Value = "ECE"
$wsPI.Cells.Item(7,5).Value = "ABE"

$wsPI.Activate()
